$wb = $excel.ActiveWorkbook

# --- Rename header labels on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add new "PO Forecast" sheet at the end ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# --- Headers ---
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# --- Data rows (ds, PO_Forecast, yhat_lower, yhat_upper) ---
$wsForecast.Range("A2").Value = 45032.99999999999
$wsForecast.Range("B2").Value = 22.0
$wsForecast.Range("C2").Value = -6.460432739584765
$wsForecast.Range("D2").Value = 51.34684346700205
$wsForecast.Range("A3").Value = 45060.99999999999
$wsForecast.Range("B3").Value = 23.0
$wsForecast.Range("C3").Value = -3.667287718968178
$wsForecast.Range("D3").Value = 56.5459692662918
$wsForecast.Range("A4").Value = 45067.99999999999
$wsForecast.Range("B4").Value = 24.0
$wsForecast.Range("C4").Value = -5.218182096691579
$wsForecast.Range("D4").Value = 54.27151408985287
$wsForecast.Range("A5").Value = 45081.99999999999
$wsForecast.Range("B5").Value = 24.0
$wsForecast.Range("C5").Value = -5.459644726952422
$wsForecast.Range("D5").Value = 54.17643544655279
$wsForecast.Range("A6").Value = 45088.99999999999
$wsForecast.Range("B6").Value = 24.0
$wsForecast.Range("C6").Value = -4.127175234943446
$wsForecast.Range("D6").Value = 53.07257938119889
$wsForecast.Range("A7").Value = 45102.99999999999
$wsForecast.Range("B7").Value = 25.0
$wsForecast.Range("C7").Value = -1.731261316804548
$wsForecast.Range("D7").Value = 52.83007066936553
$wsForecast.Range("A8").Value = 45165.99999999999
$wsForecast.Range("B8").Value = 27.0
$wsForecast.Range("C8").Value = -1.467285557032578
$wsForecast.Range("D8").Value = 56.6297319335671
$wsForecast.Range("A9").Value = 45172.99999999999
$wsForecast.Range("B9").Value = 27.0
$wsForecast.Range("C9").Value = -1.181863497090883
$wsForecast.Range("D9").Value = 54.31149191478192
$wsForecast.Range("A10").Value = 45179.99999999999
$wsForecast.Range("B10").Value = 27.0
$wsForecast.Range("C10").Value = -2.677095327327447
$wsForecast.Range("D10").Value = 53.09463825964797
$wsForecast.Range("A11").Value = 45193.99999999999
$wsForecast.Range("B11").Value = 28.0
$wsForecast.Range("C11").Value = -1.747187700678012
$wsForecast.Range("D11").Value = 56.69212514857796
$wsForecast.Range("A12").Value = 45200.99999999999
$wsForecast.Range("B12").Value = 28.0
$wsForecast.Range("C12").Value = -2.085351043283235
$wsForecast.Range("D12").Value = 55.06894847702769
$wsForecast.Range("A13").Value = 45207.99999999999
$wsForecast.Range("B13").Value = 28.0
$wsForecast.Range("C13").Value = -1.66891587910068
$wsForecast.Range("D13").Value = 55.81594426812709
$wsForecast.Range("A14").Value = 45214.99999999999
$wsForecast.Range("B14").Value = 28.0
$wsForecast.Range("C14").Value = 1.472340095114577
$wsForecast.Range("D14").Value = 59.61618076042627
$wsForecast.Range("A15").Value = 45228.99999999999
$wsForecast.Range("B15").Value = 29.0
$wsForecast.Range("C15").Value = 0.2205184410681386
$wsForecast.Range("D15").Value = 57.32370012431991
$wsForecast.Range("A16").Value = 45256.99999999999
$wsForecast.Range("B16").Value = 30.0
$wsForecast.Range("C16").Value = -2.309558585012649
$wsForecast.Range("D16").Value = 57.36767637853074
$wsForecast.Range("A17").Value = 45263.99999999999
$wsForecast.Range("B17").Value = 30.0
$wsForecast.Range("C17").Value = 0.3822471338293155
$wsForecast.Range("D17").Value = 59.67189897261738
$wsForecast.Range("A18").Value = 45277.99999999999
$wsForecast.Range("B18").Value = 30.0
$wsForecast.Range("C18").Value = 0.7167486869272605
$wsForecast.Range("D18").Value = 60.18384134220688
$wsForecast.Range("A19").Value = 45312.99999999999
$wsForecast.Range("B19").Value = 31.0
$wsForecast.Range("C19").Value = 3.076192369759479
$wsForecast.Range("D19").Value = 61.43614131522775
$wsForecast.Range("A20").Value = 45326.99999999999
$wsForecast.Range("B20").Value = 32.0
$wsForecast.Range("C20").Value = 4.940876887697329
$wsForecast.Range("D20").Value = 62.27126441546511
$wsForecast.Range("A21").Value = 45333.99999999999
$wsForecast.Range("B21").Value = 32.0
$wsForecast.Range("C21").Value = 3.064570879095774
$wsForecast.Range("D21").Value = 61.23720391659192
$wsForecast.Range("A22").Value = 45347.99999999999
$wsForecast.Range("B22").Value = 32.0
$wsForecast.Range("C22").Value = 2.28053451327204
$wsForecast.Range("D22").Value = 62.04794421211482
$wsForecast.Range("A23").Value = 45361.99999999999
$wsForecast.Range("B23").Value = 33.0
$wsForecast.Range("C23").Value = 5.169400813586539
$wsForecast.Range("D23").Value = 62.55718720932161
$wsForecast.Range("A24").Value = 45375.99999999999
$wsForecast.Range("B24").Value = 33.0
$wsForecast.Range("C24").Value = 3.597116801507329
$wsForecast.Range("D24").Value = 63.28795030761481
$wsForecast.Range("A25").Value = 45382.99999999999
$wsForecast.Range("B25").Value = 34.0
$wsForecast.Range("C25").Value = 4.827476103241962
$wsForecast.Range("D25").Value = 64.33089860505574
$wsForecast.Range("A26").Value = 45389.99999999999
$wsForecast.Range("B26").Value = 34.0
$wsForecast.Range("C26").Value = 3.702318037895006
$wsForecast.Range("D26").Value = 61.69499364697128
$wsForecast.Range("A27").Value = 45396.99999999999
$wsForecast.Range("B27").Value = 34.0
$wsForecast.Range("C27").Value = 5.590163513132005
$wsForecast.Range("D27").Value = 62.28987040948833
$wsForecast.Range("A28").Value = 45403.99999999999
$wsForecast.Range("B28").Value = 34.0
$wsForecast.Range("C28").Value = 4.646849999915115
$wsForecast.Range("D28").Value = 64.9121923873251
$wsForecast.Range("A29").Value = 45410.99999999999
$wsForecast.Range("B29").Value = 34.0
$wsForecast.Range("C29").Value = 4.048837024538701
$wsForecast.Range("D29").Value = 62.78032951049638
$wsForecast.Range("A30").Value = 45424.99999999999
$wsForecast.Range("B30").Value = 35.0
$wsForecast.Range("C30").Value = 6.350879625245477
$wsForecast.Range("D30").Value = 64.24384794720692
$wsForecast.Range("A31").Value = 45445.99999999999
$wsForecast.Range("B31").Value = 36.0
$wsForecast.Range("C31").Value = 3.139935593338463
$wsForecast.Range("D31").Value = 63.61837053933564
$wsForecast.Range("A32").Value = 45459.99999999999
$wsForecast.Range("B32").Value = 36.0
$wsForecast.Range("C32").Value = 7.236986284805374
$wsForecast.Range("D32").Value = 64.57449231987745
$wsForecast.Range("A33").Value = 45473.99999999999
$wsForecast.Range("B33").Value = 36.0
$wsForecast.Range("C33").Value = 7.194931498955371
$wsForecast.Range("D33").Value = 66.40387816440891
$wsForecast.Range("A34").Value = 45494.99999999999
$wsForecast.Range("B34").Value = 37.0
$wsForecast.Range("C34").Value = 7.439165298200277
$wsForecast.Range("D34").Value = 66.12551046761406
$wsForecast.Range("A35").Value = 45501.99999999999
$wsForecast.Range("B35").Value = 37.0
$wsForecast.Range("C35").Value = 7.26076694163665
$wsForecast.Range("D35").Value = 67.70998687578603
$wsForecast.Range("A36").Value = 45515.99999999999
$wsForecast.Range("B36").Value = 38.0
$wsForecast.Range("C36").Value = 8.481812328479364
$wsForecast.Range("D36").Value = 67.3122343243307
$wsForecast.Range("A37").Value = 45522.99999999999
$wsForecast.Range("B37").Value = 38.0
$wsForecast.Range("C37").Value = 8.48764265454588
$wsForecast.Range("D37").Value = 66.410922408438
$wsForecast.Range("A38").Value = 45529.99999999999
$wsForecast.Range("B38").Value = 38.0
$wsForecast.Range("C38").Value = 9.963351585988491
$wsForecast.Range("D38").Value = 68.01677402007128
$wsForecast.Range("A39").Value = 45536.99999999999
$wsForecast.Range("B39").Value = 38.0
$wsForecast.Range("C39").Value = 10.81029044988301
$wsForecast.Range("D39").Value = 68.0012300825985
$wsForecast.Range("A40").Value = 45543.99999999999
$wsForecast.Range("B40").Value = 39.0
$wsForecast.Range("C40").Value = 9.872920776520003
$wsForecast.Range("D40").Value = 69.92572843675265
$wsForecast.Range("A41").Value = 45550.99999999999
$wsForecast.Range("B41").Value = 39.0
$wsForecast.Range("C41").Value = 9.595180168997226
$wsForecast.Range("D41").Value = 66.62824903498026
$wsForecast.Range("A42").Value = 45557.99999999999
$wsForecast.Range("B42").Value = 39.0
$wsForecast.Range("C42").Value = 9.714350174352917
$wsForecast.Range("D42").Value = 69.29509137575678
$wsForecast.Range("A43").Value = 45564.99999999999
$wsForecast.Range("B43").Value = 39.0
$wsForecast.Range("C43").Value = 9.516756758828645
$wsForecast.Range("D43").Value = 68.66586923385098
$wsForecast.Range("A44").Value = 45571.99999999999
$wsForecast.Range("B44").Value = 40.0
$wsForecast.Range("C44").Value = 8.828951994182782
$wsForecast.Range("D44").Value = 69.82137573714354

# --- Match formatting used on the other sheets (bold/centered header, date format) ---
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A44").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Page margins matching the other sheets ---
$wsForecast.PageSetup.LeftMargin = 54
$wsForecast.PageSetup.RightMargin = 54
$wsForecast.PageSetup.TopMargin = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

$wsForecast.Range("A1").Select()
